$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29952
$ws.Range("H2").Value = 45139
$ws.Range("E3").Value = 29891
$ws.Range("F3").Value = 45078
$ws.Range("E4").Value = 29891
$ws.Range("F4").Value = 45078
$ws.Range("G4").Value = 29952
$ws.Range("H4").Value = 45139
$ws.Range("G5").Value = 29952
$ws.Range("H5").Value = 45139
$ws.Range("C6").Value = 439
$ws.Range("F6").Value = 45078
$ws.Range("G6").Value = 29952
$ws.Range("H6").Value = 45139
$ws.Range("G7").Value = 29952
$ws.Range("H7").Value = 45139
$ws.Range("D8").Value = 406
$ws.Range("E8").Value = 29891
$ws.Range("F8").Value = 45078
$ws.Range("H8").Value = 45139
$ws.Range("E9").Value = 29891
$ws.Range("F9").Value = 45078
$ws.Range("G9").Value = 29952
$ws.Range("H9").Value = 45139
$ws.Range("C10").Value = 498
$ws.Range("D10").Value = 479
$ws.Range("F10").Value = 45078
$ws.Range("H10").Value = 45139
$ws.Range("E11").Value = 29860
$ws.Range("F11").Value = 45047
$ws.Range("G11").Value = 29952
$ws.Range("H11").Value = 45139
$ws.Range("C12").Value = 367
$ws.Range("D12").Value = 348
$ws.Range("F12").Value = 45078
$ws.Range("H12").Value = 45139
$ws.Range("C13").Value = 463
$ws.Range("F13").Value = 45078
$ws.Range("G13").Value = 29952
$ws.Range("H13").Value = 45139
$ws.Range("C14").Value = 420
$ws.Range("D14").Value = 392
$ws.Range("F14").Value = 45078
$ws.Range("H14").Value = 45139
$ws.Range("C15").Value = 381
$ws.Range("F15").Value = 45108
$ws.Range("G15").Value = 29921
$ws.Range("H15").Value = 45139
$ws.Range("C16").Value = 451
$ws.Range("D16").Value = 406
$ws.Range("F16").Value = 45078
$ws.Range("H16").Value = 45139
$ws.Range("C17").Value = 366
$ws.Range("D17").Value = 390
$ws.Range("F17").Value = 45078
$ws.Range("H17").Value = 45139
$ws.Range("E18").Value = 29891
$ws.Range("F18").Value = 45078
$ws.Range("G18").Value = 29952
$ws.Range("H18").Value = 45139
$ws.Range("D19").Value = 394
$ws.Range("E19").Value = 29891
$ws.Range("F19").Value = 45078
$ws.Range("H19").Value = 45139
$ws.Range("C20").Value = 477
$ws.Range("F20").Value = 45078
$ws.Range("G20").Value = 29952
$ws.Range("H20").Value = 45139
$ws.Range("C21").Value = 306
$ws.Range("F21").Value = 45078
$ws.Range("G21").Value = 29952
$ws.Range("H21").Value = 45139
$ws.Range("C22").Value = 319
$ws.Range("D22").Value = 363
$ws.Range("F22").Value = 45078
$ws.Range("H22").Value = 45139
$ws.Range("C23").Value = 256
$ws.Range("D23").Value = 398
$ws.Range("F23").Value = 45047
$ws.Range("H23").Value = 45139
$ws.Range("C24").Value = 212
$ws.Range("D24").Value = 406
$ws.Range("F24").Value = 45108
$ws.Range("H24").Value = 45139
$ws.Range("D25").Value = 311
$ws.Range("E25").Value = 29860
$ws.Range("F25").Value = 45047
$ws.Range("H25").Value = 45139
$ws.Range("C26").Value = 329
$ws.Range("D26").Value = 309
$ws.Range("F26").Value = 45078
$ws.Range("H26").Value = 45139
$ws.Range("E27").Value = 29891
$ws.Range("F27").Value = 45078
$ws.Range("G27").Value = 29952
$ws.Range("H27").Value = 45139
$ws.Range("D28").Value = 376
$ws.Range("H28").Value = 45139
$ws.Range("C29").Value = 258
$ws.Range("D29").Value = 231
$ws.Range("F29").Value = 45078
$ws.Range("H29").Value = 45139
$ws.Range("D30").Value = 213
$ws.Range("E30").Value = 29891
$ws.Range("F30").Value = 45078
$ws.Range("H30").Value = 45139
$ws.Range("C31").Value = 390
$ws.Range("F31").Value = 45078
$ws.Range("G31").Value = 29952
$ws.Range("H31").Value = 45139
$ws.Range("E32").Value = 29891
$ws.Range("F32").Value = 45078
$ws.Range("G32").Value = 29952
$ws.Range("H32").Value = 45139
$ws.Range("D33").Value = 394
$ws.Range("H33").Value = 45139
$ws.Range("C34").Value = 198
$ws.Range("D34").Value = 313
$ws.Range("F34").Value = 45078
$ws.Range("H34").Value = 45139
$ws.Range("C35").Value = 400
$ws.Range("D35").Value = 313
$ws.Range("F35").Value = 45047
$ws.Range("H35").Value = 45139
$ws.Range("D36").Value = 406
$ws.Range("H36").Value = 45139
$ws.Range("C37").Value = 462
$ws.Range("D37").Value = 313
$ws.Range("F37").Value = 45078
$ws.Range("H37").Value = 45139
$ws.Range("C38").Value = 355
$ws.Range("D38").Value = 363
$ws.Range("F38").Value = 45078
$ws.Range("H38").Value = 45139
$ws.Range("C39").Value = 223
$ws.Range("D39").Value = 219
$ws.Range("F39").Value = 45078
$ws.Range("H39").Value = 45139
$ws.Range("C40").Value = 282
$ws.Range("D40").Value = 310
$ws.Range("F40").Value = 45078
$ws.Range("H40").Value = 45139
$ws.Range("C41").Value = 390
$ws.Range("D41").Value = 312
$ws.Range("F41").Value = 45078
$ws.Range("H41").Value = 45139
$ws.Range("D42").Value = 214
$ws.Range("H42").Value = 45139
$ws.Range("C43").Value = 469
$ws.Range("D43").Value = 313
$ws.Range("F43").Value = 45078
$ws.Range("H43").Value = 45139
$ws.Range("D44").Value = 300
$ws.Range("H44").Value = 45139
$ws.Range("D45").Value = 313
$ws.Range("H45").Value = 45139
$ws.Range("C46").Value = 330
$ws.Range("D46").Value = 294
$ws.Range("F46").Value = 45078
$ws.Range("H46").Value = 45139
$ws.Range("C47").Value = 331
$ws.Range("D47").Value = 256
$ws.Range("F47").Value = 45078
$ws.Range("H47").Value = 45139
$ws.Range("C48").Value = 354
$ws.Range("D48").Value = 311
$ws.Range("F48").Value = 45078
$ws.Range("H48").Value = 45139
$ws.Range("C49").Value = 293
$ws.Range("D49").Value = 309
$ws.Range("F49").Value = 45047
$ws.Range("H49").Value = 45139
$ws.Range("C50").Value = 355
$ws.Range("D50").Value = 234
$ws.Range("F50").Value = 45078
$ws.Range("H50").Value = 45139
$ws.Range("D51").Value = 313
$ws.Range("E51").Value = 29707
$ws.Range("F51").Value = 45047
$ws.Range("H51").Value = 45139
$ws.Range("D52").Value = 311
$ws.Range("H52").Value = 45139

Write-Host "Updated 182 cells"
